$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ResourceNode")

# Add new row 6 - a test resource node entry (copy of JSM row's drop-table refs)
$ws.Range("B6").Value = "NOD_HER_TEST_001"
$ws.Range("C6").Value = "테스트"
$ws.Range("D6").Value = "Test_Node"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 30
$ws.Range("G6").Value = 4
$ws.Range("H6").Value = $true
$ws.Range("I6").Value = $false
$ws.Range("J6").Value = "JSM_INT"
$ws.Range("K6").Value = "JSM_DST"
$ws.Range("L6").Value = "JSM_FIN"
$ws.Range("M6").Value = "temporary"
$ws.Range("N6").Value = "JSM_CON"
$ws.Range("O6").Value = "쟈스민 기본 채집물"

# Move the active selection on sheet1 to A4
$ws.Range("A4").Select()

# Move the active selection on sheet2 to B5
$ws2 = $wb.Worksheets.Item("ResourceNode_DropTable")
$ws2.Range("B5").Select()
